# Reading of the input file was done.
#
# - Rename "Reels 95.5 RTP" -> "Base Reels 95.5 RTP"
# - Add a new, empty "Free Reels 95.5 RTP" sheet right after it
# - Zoom the new sheet to 130%
# - Leave "Base Reels 95.5 RTP" as the selected / active sheet

$wb = $excel.ActiveWorkbook

$baseReels = $wb.Worksheets.Item("Reels 95.5 RTP")
$baseReels.Name = "Base Reels 95.5 RTP"

$freeReels = $wb.Worksheets.Add($null, $baseReels)
$freeReels.Name = "Free Reels 95.5 RTP"

$freeReels.Activate()
$excel.ActiveWindow.Zoom = 130

$baseReels.Activate()
$baseReels.Select()
